$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly fruit/vegetable data refresh: the "Níspero" lot records for the
# Macroferia Regional de Talca effectively swap between the two reporting
# dates (rows 2-3 <-> rows 4-5), updating date, variety, quality, volume,
# min/max/average price, commercialization unit and $/Kg price.

# Row 2 (was 2021-11-05 Californiana(o)/Primera) becomes the
# 2022-12-07 Golden Nugget/Especial record.
$ws.Range("D2").Value2 = 44902
$ws.Range("K2").Value2 = "Golden Nugget"
$ws.Range("L2").Value2 = "Especial"
$ws.Range("M2").Value2 = 60
$ws.Range("N2").Value2 = 15000
$ws.Range("O2").Value2 = 15000
$ws.Range("P2").Value2 = 15000
$ws.Range("Q2").Value2 = "$/caja 10 kilos"
$ws.Range("S2").Value2 = 1500

# Row 3 (was 2021-11-05 Golden Nugget/Primera) becomes the
# 2022-12-07 Golden Nugget/Primera record with updated prices.
$ws.Range("D3").Value2 = 44902
$ws.Range("M3").Value2 = 70
$ws.Range("N3").Value2 = 13000
$ws.Range("O3").Value2 = 13000
$ws.Range("P3").Value2 = 13000
$ws.Range("Q3").Value2 = "$/caja 10 kilos"
$ws.Range("S3").Value2 = 1300

# Row 4 (was 2022-12-07 Golden Nugget/Especial) becomes the
# 2021-11-05 Californiana(o)/Primera record.
$ws.Range("D4").Value2 = 44505
$ws.Range("K4").Value2 = "Californiana(o)"
$ws.Range("L4").Value2 = "Primera"
$ws.Range("M4").Value2 = 100
$ws.Range("Q4").Value2 = "$/bandeja 10 kilos"

# Row 5 (was 2022-12-07 Golden Nugget/Primera) becomes the
# 2021-11-05 Golden Nugget/Primera record with updated prices.
$ws.Range("D5").Value2 = 44505
$ws.Range("M5").Value2 = 50
$ws.Range("N5").Value2 = 15000
$ws.Range("O5").Value2 = 15000
$ws.Range("P5").Value2 = 15000
$ws.Range("Q5").Value2 = "$/bandeja 10 kilos"
$ws.Range("S5").Value2 = 1500
